# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("N6").Value = 15
$ws.Range("O6").Value = 1.18
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.6
$ws.Range("R6").Value = 2.3
$ws.Range("U6").Value = 1.75
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 8
$ws.Range("Z6").Value = 11
$ws.Range("AA6").Value = 12
$ws.Range("AC6").Value = 15
$ws.Range("AD6").Value = 8.5
$ws.Range("AG6").Value = 201
$ws.Range("AK6").Value = 67
$ws.Range("AN6").Value = 3.6
$ws.Range("AP6").Value = 17
$ws.Range("AV6").Value = 51
$ws.Range("AW6").Value = 501
$ws.Range("BA6").Value = 101
$ws.Range("BC6").Value = 201

# Row 7
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 2.1

# Row 8
$ws.Range("Q8").Value = 1.5
$ws.Range("R8").Value = 2.5
$ws.Range("U8").Value = 1.57
$ws.Range("V8").Value = 2.25
$ws.Range("AW8").Value = 351

# Row 10
$ws.Range("N10").Value = 15
$ws.Range("O10").Value = 1.18
$ws.Range("P10").Value = 4.5
$ws.Range("Q10").Value = 1.6
$ws.Range("R10").Value = 2.3

# Row 11
$ws.Range("G11").Value = 2.3
$ws.Range("I11").Value = 2.7
$ws.Range("X11").Value = 17
$ws.Range("AA11").Value = 15
$ws.Range("AC11").Value = 23
$ws.Range("AZ11").Value = 17
$ws.Range("BC11").Value = 81

# Row 24
$ws.Range("N24").Value = 9
